$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ApplicationLogin")

# E8: Test_Data was "NA" -> becomes numeric 10
$ws.Range("E8").Value = 10

# E9: Test_Data trailing-space trimmed: "Dashboard « Stock Accounting " -> "Dashboard « Stock Accounting"
$ws.Range("E9").Value = "Dashboard « Stock Accounting"

# Selection moved from B11 to E9
$ws.Range("E9").Select()

# Rows 12-17 lose their borders (style changes from bordered "s=10" to borderless "s=8")
$ws.Range("A12:F17").Borders.LineStyle = -4142

